$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix cell G52: cost value 4 -> 1 ---
$ws.Range("G52").Value = 1

# --- Remove the stray explicit "applyFont" style from V40 (now unused cellXf) ---
$ws.Range("V40").ClearFormats()

# --- New card rows 69-87 (ids/group already existed; fill remaining columns) ---
# NOTE: two "ability" cells (O71, O86) are intentionally written last, below,
# to reproduce the exact shared-string insertion order of the source workbook.
# row 69
$ws.Range("C69").Value = "Kawakawa"
$ws.Range("D69").Value = "Euthynnus affinis"
$ws.Range("E69").Value = 1
$ws.Range("F69").Value = 2
$ws.Range("I69").Value = 1
$ws.Range("L69").Value = 6
$ws.Range("M69").Value = 100
$ws.Range("N69").Value = "WhenPlayed"
$ws.Range("O69").Value = "[FishEgg][ArrowDown][FlipperBlue]"
$ws.Range("V69").Value = "Adults live in open waters but remain close to the shoreline. Their young may enter bays and harbors."

# row 70
$ws.Range("C70").Value = "Largetooth Flounder"
$ws.Range("D70").Value = "Pseudorhombus arsius"
$ws.Range("E70").Value = 3
$ws.Range("I70").Value = 1
$ws.Range("L70").Value = 5
$ws.Range("M70").Value = 45
$ws.Range("N70").Value = "GameEnd"
$ws.Range("O70").Value = "[FishEgg][ArrowDown][FlipperGreen]"
$ws.Range("Q70").Value = 1
$ws.Range("V70").Value = "The largetooth flounder is a species of left- eyed flatfish with both eyes of the adult on the left side of its body."

# row 71
$ws.Range("C71").Value = "Leafy Seadragon"
$ws.Range("D71").Value = "Phycodurus eques"
$ws.Range("F71").Value = 1
$ws.Range("I71").Value = 1
$ws.Range("L71").Value = 2
$ws.Range("M71").Value = 35
$ws.Range("N71").Value = "GameEnd"
$ws.Range("Q71").Value = 1
$ws.Range("V71").Value = "Lobes of skin growing all over it provide camouflage, giving it the appearance of seaweed, even as it swims."

# row 73
$ws.Range("C73").Value = "Live Sharksucker"
$ws.Range("D73").Value = "Echeneis naucrates"
$ws.Range("F73").Value = 2
$ws.Range("I73").Value = 1
$ws.Range("L73").Value = 4
$ws.Range("M73").Value = 110
$ws.Range("N73").Value = "WhenPlayed"
$ws.Range("O73").Value = "[ConsumeFish1][ConsumeFish1]"
$ws.Range("V73").Value = "This fish can attach itself to any large host, even scuba divers, using the oval-shaped sucking disc on its head."

# row 74
$ws.Range("C74").Value = "Long-Snouted Lancetfish"
$ws.Range("D74").Value = "Alepisaurus ferox"
$ws.Range("H74").Value = 1
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1
$ws.Range("K74").Value = 1
$ws.Range("L74").Value = 3
$ws.Range("M74").Value = 200
$ws.Range("N74").Value = "WhenPlayed"
$ws.Range("O74").Value = "[DrawCard][DrawCard][DrawCard][SchoolFeederMove]"
$ws.Range("S74").Value = 1
$ws.Range("V74").Value = "Voracious eaters, many new species of fish and mollusks have been discovered in the contents of their stomachs."

# row 75
$ws.Range("C75").Value = "Longspine Porcupinefish"
$ws.Range("D75").Value = "Diodon holocanthus"
$ws.Range("E75").Value = 2
$ws.Range("F75").Value = 1
$ws.Range("I75").Value = 1
$ws.Range("L75").Value = 5
$ws.Range("M75").Value = 50
$ws.Range("N75").Value = "WhenPlayed"
$ws.Range("O75").Value = "[FishEgg][ArrowDown][FishLengthMedium] on each"
$ws.Range("T75").Value = 1
$ws.Range("V75").Value = "It uses its beak and the plates on the roof of its mouth to crush prey that would otherwise be indigestible."

# row 76
$ws.Range("C76").Value = "Mahi-Mahi"
$ws.Range("D76").Value = "Coryphaena hippurus"
$ws.Range("E76").Value = 1
$ws.Range("F76").Value = 2
$ws.Range("I76").Value = 1
$ws.Range("L76").Value = 7
$ws.Range("M76").Value = 210
$ws.Range("N76").Value = "IfActivated"
$ws.Range("O76").Value = "[FishEgg]"
$ws.Range("S76").Value = 1
$ws.Range("V76").Value = "Its common name means “strong-strong” in Hawaiian. Flyingfishes make up nearly 25% of the diet of adults."

# row 77
$ws.Range("C77").Value = "Maletese Ray"
$ws.Range("D77").Value = "Leucoraja melitensis"
$ws.Range("F77").Value = 2
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 1
$ws.Range("L77").Value = 3
$ws.Range("M77").Value = 40
$ws.Range("N77").Value = "WhenPlayed"
$ws.Range("O77").Value = "[FishEgg][ArrowDown][FishLengthSmall] on each"
$ws.Range("Q77").Value = 1
$ws.Range("V77").Value = "Critically endangered due to habitat loss, the Maltese ray faces a high risk of extinction in the wild."

# row 78
$ws.Range("C78").Value = "Mariana Snailfish"
$ws.Range("D78").Value = "Pseudoliparis swirei"
$ws.Range("E78").Value = 2
$ws.Range("K78").Value = 2
$ws.Range("L78").Value = 1
$ws.Range("M78").Value = 28
$ws.Range("N78").Value = "GameEnd"
$ws.Range("O78").Value = "[FishEgg][ArrowDown][PlayFishBottomRow] on each"
$ws.Range("V78").Value = "It has been observed at a depth of 8,178 m, making it one of the deepest dwelling fishes on the planet."

# row 79
$ws.Range("C79").Value = "Mediterranean Damselfish"
$ws.Range("D79").Value = "Chromis chromis"
$ws.Range("G79").Value = 1
$ws.Range("I79").Value = 1
$ws.Range("L79").Value = 4
$ws.Range("M79").Value = 25
$ws.Range("N79").Value = "IfActivated"
$ws.Range("O79").Value = "[FishEgg]"
$ws.Range("U79").Value = "purple"
$ws.Range("V79").Value = "Though normally a peaceful species, males aggressively guard the eggs while they develop."

# row 80
$ws.Range("C80").Value = "Mediterranean Parrotfish"
$ws.Range("D80").Value = "Sparisoma cretense"
$ws.Range("E80").Value = 2
$ws.Range("F80").Value = 1
$ws.Range("I80").Value = 1
$ws.Range("L80").Value = 5
$ws.Range("M80").Value = 50
$ws.Range("N80").Value = "IfActivated"
$ws.Range("O80").Value = "[FishHatch]"
$ws.Range("U80").Value = "green"
$ws.Range("V80").Value = "It starts life as female, then changes to male in adulthood, turning from bright red and yellow to gray and blue."

# row 81
$ws.Range("C81").Value = "Midwater Scorpionfish"
$ws.Range("D81").Value = "Ectreposebastes imus"
$ws.Range("E81").Value = 3
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 1
$ws.Range("K81").Value = 1
$ws.Range("L81").Value = 3
$ws.Range("M81").Value = 18
$ws.Range("N81").Value = "IfActivated"
$ws.Range("O81").Value = "[SchoolFeederMove]"
$ws.Range("T81").Value = 1
$ws.Range("V81").Value = "This fish has often been seen suspending itself vertically in the water, with its mouth pointing up towards the surface."

# row 82
$ws.Range("C82").Value = "Ocean sunfish"
$ws.Range("D82").Value = "Mola mola"
$ws.Range("E82").Value = 1
$ws.Range("F82").Value = 2
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 1
$ws.Range("L82").Value = 4
$ws.Range("M82").Value = 330
$ws.Range("N82").Value = "GameEnd"
$ws.Range("O82").Value = "[FishEgg][ArrowDown][FishLengthLarge] on each"
$ws.Range("V82").Value = "The largest known bony fish, it basks on the surface to thermally recharge after feeding in deeper, colder waters."

# row 83
$ws.Range("C83").Value = "Oceanic Puffer"
$ws.Range("D83").Value = "Lagocephalus lagocephalus"
$ws.Range("E83").Value = 2
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 1
$ws.Range("L83").Value = 5
$ws.Range("M83").Value = 60
$ws.Range("N83").Value = "IfActivated"
$ws.Range("O83").Value = "(all players) [SchoolFeederMove][AllPlayers]"
$ws.Range("T83").Value = 1
$ws.Range("V83").Value = "Its Latin name means “rabbit head,” and, like many puffers, it contains neurotoxins that may be fatal to humans if eaten."

# row 84
$ws.Range("C84").Value = "Pacific Sardine"
$ws.Range("D84").Value = "Sardinops sagax"
$ws.Range("E84").Value = 3
$ws.Range("I84").Value = 1
$ws.Range("L84").Value = 4
$ws.Range("M84").Value = 40
$ws.Range("N84").Value = "WhenPlayed"
$ws.Range("O84").Value = "[FishEgg][FishEgg][FishHatch][FishHatch]"
$ws.Range("V84").Value = "Forming schools of up to 10 million fish, it is, at times, the most abundant species in the California Current."

# row 85
$ws.Range("C85").Value = "Pacific White Skate"
$ws.Range("D85").Value = "Bathyraja spinosissima"
$ws.Range("F85").Value = 2
$ws.Range("J85").Value = 1
$ws.Range("K85").Value = 1
$ws.Range("L85").Value = 3
$ws.Range("M85").Value = 150
$ws.Range("N85").Value = "WhenPlayed"
$ws.Range("O85").Value = "[FishHatch][YoungFish]"
$ws.Range("V85").Value = "Most skates only have dermal denticles on their dorsal side, but this one has tiny spines on its belly as well."

# row 86
$ws.Range("C86").Value = "Pale Chimaera"
$ws.Range("D86").Value = "Hydrolagus pallidus"
$ws.Range("E86").Value = 1
$ws.Range("F86").Value = 1
$ws.Range("J86").Value = 1
$ws.Range("K86").Value = 1
$ws.Range("L86").Value = 2
$ws.Range("M86").Value = 137
$ws.Range("N86").Value = "GameEnd"
$ws.Range("V86").Value = "Dubbed “ghost sharks,” chimaeras “flap” their large pectoral fins to propel themselves forward (unlike sharks)."

# row 87
$ws.Range("C87").Value = "Paintspotted Moray"
$ws.Range("D87").Value = "Gymnothorax pictus"
$ws.Range("F87").Value = 2
$ws.Range("I87").Value = 1
$ws.Range("L87").Value = 4
$ws.Range("M87").Value = 140
$ws.Range("N87").Value = "WhenPlayed"
$ws.Range("O87").Value = "[Discard][Discard][Discard]"
$ws.Range("S87").Value = 1
$ws.Range("V87").Value = "It inhabits reef flats and rocky intertidal shorelines, and it will sometimes leave the water in pursuit of prey."

# --- Deferred ability-text cells (written last originally) ---
$ws.Range("O71").Value = "3 [Wave] if [SchoolFish] on this fish"
$ws.Range("O86").Value = "10 [Wave] if [FishEgg] + [YoungFish] + [SchoolFish] on this fish"

# --- Update the saved selection / scroll position to match the editor state at save time ---
$ws.Range("O83").Select()
